$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric metrics.

$data = @{
    2 = @("FAPs", "Wnt1", "Ror2", "ECs", 1, 0.3333333333333333, 0.01948966666666667, 0.058469, 0.07096062449330311, 0.07096062449330311, 1, 0.3333333333333333, 0.047572, 0.142716, 0.006780712238476238, 0.00730352994735907, 0.0009271624226666667, 0.008344461804000001, 0.0004811635749516571, 0.0005182630460701408)
    3 = @("FAPs", "Wnt1", "Ror2", "FAPs", 1, 0.3333333333333333, 0.01948966666666667, 0.058469, 0.07096062449330311, 0.07096062449330311, 3, 1, 5.461551666666666, 16.384655, 0.7784665397132128, 0.8384891565742209, 0.1064438214661111, 0.9579943931949999, 0.05524047180519032, 0.05949971418136972)
    4 = @("FAPs", "Wnt1", "Ror2", "sCs", 1, 0.3333333333333333, 0.01948966666666667, 0.058469, 0.07096062449330311, 0.07096062449330311, 2, 1, 1.5066585, 3.013317, 0.214752748048311, 0.15420731347842, 0.0293642719455, 0.176185631673, 0.01523898911316113, 0.01094264726586324)
    5 = @("ECs", "Wnt1", "Ror2", "FAPs", 3, 1, 0.255165, 0.7654949999999999, 0.9290393755066968, 0.9290393755066969, 1, 0.3333333333333333, 0.047572, 0.142716, 0.006780712238476238, 0.00730352994735907, 0.01213870938, 0.10924838442, 0.00629954866352458, 0.00678526690128893)
    6 = @("ECs", "Wnt1", "Ror2", "FAPs", 3, 1, 0.255165, 0.7654949999999999, 0.9290393755066968, 0.9290393755066969, 3, 1, 5.461551666666666, 16.384655, 0.7784665397132128, 0.8384891565742209, 1.393596831025, 12.542371479225, 0.7232260679080224, 0.7789894423928512)
    7 = @("ECs", "Wnt1", "Ror2", "sCs", 3, 1, 0.255165, 0.7654949999999999, 0.9290393755066968, 0.9290393755066969, 2, 1, 1.5066585, 3.013317, 0.214752748048311, 0.15420731347842, 0.3844465161524999, 2.306679096915, 0.1995137589351498, 0.1432646662125568)
}

foreach ($r in 2..7) {
    $vals = $data[$r]
    for ($c = 1; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $vals[$c - 1]
    }
}
